$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.5
$ws.Range("G2").Value = 2.52
$ws.Range("H2").Value = 3.4
$ws.Range("I2").Value = 3.45
$ws.Range("J2").Value = 3.2
$ws.Range("L2").Value = 1.49
$ws.Range("Q2").Value = 2.18
$ws.Range("S2").Value = 4
$ws.Range("T2").Value = 1.84
$ws.Range("V2").Value = 1.38
$ws.Range("W2").Value = 1.62
$ws.Range("X2").Value = 11.5
$ws.Range("Y2").Value = 12
$ws.Range("Z2").Value = 22
$ws.Range("AA2").Value = 60
$ws.Range("AB2").Value = 9.800000000000001
$ws.Range("AD2").Value = 14.5
$ws.Range("AE2").Value = 42
$ws.Range("AF2").Value = 15.5
$ws.Range("AH2").Value = 18
$ws.Range("AJ2").Value = 34
$ws.Range("AK2").Value = 30
$ws.Range("AL2").Value = 46
$ws.Range("F3").Value = 1.71
$ws.Range("G3").Value = 1.76
$ws.Range("H3").Value = 6.2
$ws.Range("I3").Value = 6.6
$ws.Range("J3").Value = 3.65
$ws.Range("K3").Value = 3.9
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 2.94
$ws.Range("O3").Value = 1.48
$ws.Range("P3").Value = 1.63
$ws.Range("Q3").Value = 2.46
$ws.Range("T3").Value = 2.28
$ws.Range("U3").Value = 1.7
$ws.Range("V3").Value = 1.14
$ws.Range("W3").Value = 2.06
$ws.Range("X3").Value = 10.5
$ws.Range("Y3").Value = 16.5
$ws.Range("Z3").Value = 48
$ws.Range("AA3").Value = 220
$ws.Range("AB3").Value = 6.6
$ws.Range("AC3").Value = 8.6
$ws.Range("AD3").Value = 27
$ws.Range("AE3").Value = 130
$ws.Range("AF3").Value = 8.800000000000001
$ws.Range("AH3").Value = 30
$ws.Range("AI3").Value = 150
$ws.Range("AJ3").Value = 18
$ws.Range("AL3").Value = 60
$ws.Range("AM3").Value = 240
$ws.Range("AN3").Value = 16
$ws.Range("AO3").Value = 220
$ws.Range("G4").Value = 2.96
$ws.Range("I4").Value = 2.92
$ws.Range("J4").Value = 3.15
$ws.Range("N4").Value = 3.05
$ws.Range("V4").Value = 1.52
$ws.Range("X4").Value = 10
$ws.Range("Y4").Value = 9.800000000000001
$ws.Range("AD4").Value = 13.5
$ws.Range("AF4").Value = 18.5
$ws.Range("AN4").Value = 1000
$ws.Range("AO4").Value = 48
$ws.Range("F5").Value = 1.64
$ws.Range("G5").Value = 1.69
$ws.Range("I5").Value = 6.8
$ws.Range("J5").Value = 3.9
$ws.Range("K5").Value = 4.2
$ws.Range("N5").Value = 3.1
$ws.Range("O5").Value = 1.42
$ws.Range("P5").Value = 1.75
$ws.Range("Q5").Value = 2.26
$ws.Range("R5").Value = 1.26
$ws.Range("T5").Value = 2.2
$ws.Range("W5").Value = 2.44
$ws.Range("Y5").Value = 17
$ws.Range("AA5").Value = 250
$ws.Range("AC5").Value = 9.800000000000001
$ws.Range("AF5").Value = 9.199999999999999
$ws.Range("AH5").Value = 30
$ws.Range("AN5").Value = 16
$ws.Range("F6").Value = 1.85
$ws.Range("G6").Value = 1.86
$ws.Range("H6").Value = 5.3
$ws.Range("I6").Value = 5.9
$ws.Range("J6").Value = 3.5
$ws.Range("K6").Value = 3.65
$ws.Range("M6").Value = 1.09
$ws.Range("N6").Value = 3.3
$ws.Range("O6").Value = 1.38
$ws.Range("P6").Value = 1.77
$ws.Range("Q6").Value = 2.12
$ws.Range("R6").Value = 1.29
$ws.Range("T6").Value = 2
$ws.Range("U6").Value = 1.8
$ws.Range("V6").Value = 1.21
$ws.Range("W6").Value = 2.16
$ws.Range("X6").Value = 12
$ws.Range("Y6").Value = 17
$ws.Range("AA6").Value = 170
$ws.Range("AB6").Value = 7.8
$ws.Range("AD6").Value = 22
$ws.Range("AF6").Value = 10.5
$ws.Range("AH6").Value = 21
$ws.Range("AI6").Value = 1000
$ws.Range("AJ6").Value = 19.5
$ws.Range("AK6").Value = 25
$ws.Range("AL6").Value = 48
$ws.Range("AM6").Value = 160
$ws.Range("AN6").Value = 16
$ws.Range("F7").Value = 1.08
$ws.Range("H7").Value = 34
$ws.Range("I7").Value = 1000
$ws.Range("J7").Value = 15.5
$ws.Range("L7").Value = 1.2
$ws.Range("M7").Value = 1.01
$ws.Range("N7").Value = 7.6
$ws.Range("O7").Value = 1.11
$ws.Range("P7").Value = 3.35
$ws.Range("Q7").Value = 1.36
$ws.Range("R7").Value = 1.98
$ws.Range("S7").Value = 1.94
$ws.Range("U7").Value = 1.36
$ws.Range("V7").Value = 1.01
$ws.Range("W7").Value = 11
$ws.Range("AB7").Value = 980
$ws.Range("AJ7").Value = 8.6
$ws.Range("F8").Value = 1.65
$ws.Range("G8").Value = 1.69
$ws.Range("I8").Value = 7.4
$ws.Range("J8").Value = 3.85
$ws.Range("K8").Value = 4.1
$ws.Range("N8").Value = 2.92
$ws.Range("P8").Value = 1.65
$ws.Range("Q8").Value = 2.42
$ws.Range("T8").Value = 2.32
$ws.Range("U8").Value = 1.64
$ws.Range("W8").Value = 2.44
$ws.Range("X8").Value = 10.5
$ws.Range("Y8").Value = 17.5
$ws.Range("AA8").Value = 330
$ws.Range("AC8").Value = 9
$ws.Range("AD8").Value = 29
$ws.Range("AE8").Value = 140
$ws.Range("AF8").Value = 8.199999999999999
$ws.Range("AI8").Value = 160
$ws.Range("AJ8").Value = 16
$ws.Range("AN8").Value = 15
$ws.Range("F9").Value = 1.52
$ws.Range("I9").Value = 9
$ws.Range("K9").Value = 4.4
$ws.Range("V9").Value = 1.13
$ws.Range("AK9").Value = 1000
$ws.Range("F10").Value = 1.4
$ws.Range("G10").Value = 1.42
$ws.Range("H10").Value = 10
$ws.Range("I10").Value = 12.5
$ws.Range("M10").Value = 1.06
$ws.Range("P10").Value = 2.08
$ws.Range("Q10").Value = 1.82
$ws.Range("R10").Value = 1.43
$ws.Range("S10").Value = 3
$ws.Range("T10").Value = 2.12
$ws.Range("U10").Value = 1.72
$ws.Range("V10").Value = 1.09
$ws.Range("W10").Value = 3.35
$ws.Range("X10").Value = 18
$ws.Range("Y10").Value = 32
$ws.Range("Z10").Value = 100
$ws.Range("AA10").Value = 470
$ws.Range("AB10").Value = 8.199999999999999
$ws.Range("AD10").Value = 42
$ws.Range("AE10").Value = 200
$ws.Range("AF10").Value = 8
$ws.Range("AG10").Value = 10.5
$ws.Range("AH10").Value = 32
$ws.Range("AI10").Value = 170
$ws.Range("AJ10").Value = 11.5
$ws.Range("AK10").Value = 16
$ws.Range("AL10").Value = 42
$ws.Range("AN10").Value = 7.2
$ws.Range("AO10").Value = 260
$ws.Range("N11").Value = 3.8
$ws.Range("O11").Value = 1.21
$ws.Range("U11").Value = 2.26
$ws.Range("F12").Value = 2.6
$ws.Range("G12").Value = 2.82
$ws.Range("H12").Value = 2.92
$ws.Range("I12").Value = 3.2
$ws.Range("J12").Value = 3.2
$ws.Range("K12").Value = 3.55
$ws.Range("O12").Value = 1.45
$ws.Range("P12").Value = 1.64
$ws.Range("Q12").Value = 2.34
$ws.Range("T12").Value = 1.93
$ws.Range("V12").Value = 1.46
$ws.Range("W12").Value = 1.55
$ws.Range("X12").Value = 13.5
$ws.Range("Y12").Value = 14
$ws.Range("Z12").Value = 21
$ws.Range("AA12").Value = 60
$ws.Range("AB12").Value = 10
$ws.Range("AD12").Value = 1000
$ws.Range("AE12").Value = 1000
$ws.Range("AF12").Value = 17
$ws.Range("AG12").Value = 12.5
$ws.Range("AJ12").Value = 42
$ws.Range("AK12").Value = 38
$ws.Range("AL12").Value = 60
$ws.Range("AN12").Value = 38
$ws.Range("AO12").Value = 50
